$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 63666.668
$ws.Range("I11").Value = 63666.668
$ws.Range("K11").Value = 63666.668
$ws.Range("M11").Value = -63526.668
$ws.Range("H33").Value = 107681.29
$ws.Range("I33").Value = 150313.4
$ws.Range("K33").Value = 150313.4
$ws.Range("M33").Value = -150084.4
$ws.Range("H43").Value = 1869.8
$ws.Range("J43").Value = 2924.75
$ws.Range("L43").Value = 2924.75
$ws.Range("N43").Value = -3062.75
$ws.Range("H100").Value = 6564.846
$ws.Range("I100").Value = 9892.166999999999
$ws.Range("K100").Value = 9892.166999999999
$ws.Range("M100").Value = -9351.166999999999
$ws.Range("H132").Value = 2987.2454
$ws.Range("I132").Value = 2712.2744
$ws.Range("K132").Value = 8136.823199999999
$ws.Range("M132").Value = -5606.823199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3500.7917
$ws.Range("I61").Value = 1879.9333
$ws.Range("J61").Value = 6202.222
$ws.Range("K61").Value = 1879.9333
$ws.Range("L61").Value = 6202.222
$ws.Range("M61").Value = -1667.9333
$ws.Range("N61").Value = -6626.222
$ws.Range("H74").Value = 50004230
$ws.Range("I74").Value = 83335640
$ws.Range("J74").Value = 7124.375
$ws.Range("K74").Value = 83335640
$ws.Range("L74").Value = 7124.375
$ws.Range("M74").Value = -83334766
$ws.Range("N74").Value = -8872.375
$ws.Range("H77").Value = 50004230
$ws.Range("I77").Value = 83335640
$ws.Range("J77").Value = 7124.375
$ws.Range("K77").Value = 416678200
$ws.Range("L77").Value = 35621.875
$ws.Range("M77").Value = -416673832
$ws.Range("N77").Value = -44357.875
$ws.Range("H97").Value = 1068.8096
$ws.Range("I97").Value = 1132.1177
$ws.Range("K97").Value = 1132.1177
$ws.Range("M97").Value = -636.1177
$ws.Range("H102").Value = 1957088.5
$ws.Range("I102").Value = 2444500.2
$ws.Range("K102").Value = 2444500.2
$ws.Range("M102").Value = -2442878.2
$ws.Range("H132").Value = 7720.1113
$ws.Range("I132").Value = 10012
$ws.Range("J132").Value = 7433.625
$ws.Range("K132").Value = 30036
$ws.Range("L132").Value = 22300.875
$ws.Range("M132").Value = -27506
$ws.Range("N132").Value = -27360.875
$ws.Range("H136").Value = 3500.7917
$ws.Range("I136").Value = 1879.9333
$ws.Range("J136").Value = 6202.222
$ws.Range("K136").Value = 5639.7999
$ws.Range("L136").Value = 18606.666
$ws.Range("M136").Value = -3089.7999
$ws.Range("N136").Value = -23706.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 25003218
$ws.Range("I94").Value = 41669030
$ws.Range("J94").Value = 4496
$ws.Range("K94").Value = 41669030
$ws.Range("L94").Value = 4496
$ws.Range("M94").Value = -41668579
$ws.Range("N94").Value = -5398
$ws.Range("H99").Value = 1927.375
$ws.Range("I99").Value = 1454
$ws.Range("J99").Value = 1995
$ws.Range("K99").Value = 1454
$ws.Range("L99").Value = 1995
$ws.Range("M99").Value = 44
$ws.Range("N99").Value = -4991
$ws.Range("H134").Value = 2952.5
$ws.Range("I134").Value = 2428.3225
$ws.Range("K134").Value = 7284.967500000001
$ws.Range("M134").Value = -4749.967500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3287.3235
$ws.Range("J31").Value = 3500.4644
$ws.Range("L31").Value = 3500.4644
$ws.Range("N31").Value = -4090.4644
$ws.Range("H34").Value = 3287.3235
$ws.Range("J34").Value = 3500.4644
$ws.Range("L34").Value = 3500.4644
$ws.Range("N34").Value = -3904.4644
$ws.Range("H60").Value = 8250
$ws.Range("I60").Value = 8250
$ws.Range("K60").Value = 8250
$ws.Range("M60").Value = -7739
$ws.Range("H62").Value = 133999.6
$ws.Range("I62").Value = 20000
$ws.Range("J62").Value = 162499.5
$ws.Range("K62").Value = 20000
$ws.Range("L62").Value = 162499.5
$ws.Range("M62").Value = -19376
$ws.Range("N62").Value = -163747.5
$ws.Range("H65").Value = 133999.6
$ws.Range("I65").Value = 20000
$ws.Range("J65").Value = 162499.5
$ws.Range("K65").Value = 100000
$ws.Range("L65").Value = 812497.5
$ws.Range("M65").Value = -96880
$ws.Range("N65").Value = -818737.5
$ws.Range("H107").Value = 724.75
$ws.Range("I107").Value = 819
$ws.Range("J107").Value = 442
$ws.Range("K107").Value = 819
$ws.Range("L107").Value = 442
$ws.Range("M107").Value = 1101
$ws.Range("N107").Value = -4282
$ws.Range("H135").Value = 63298.43
$ws.Range("J135").Value = 63298.43
$ws.Range("L135").Value = 63298.43
$ws.Range("N135").Value = -73438.42999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 385.96667
$ws.Range("I2").Value = 20.411764
$ws.Range("J2").Value = 864
$ws.Range("K2").Value = 122.470584
$ws.Range("L2").Value = 5184
$ws.Range("M2").Value = -9.470584000000002
$ws.Range("N2").Value = -5410
$ws.Range("H4").Value = 33845196
$ws.Range("I4").Value = 35383600
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 106150800
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = -106150688
$ws.Range("N4").Value = -1124
$ws.Range("H6").Value = 54.125
$ws.Range("I6").Value = 39.833332
$ws.Range("J6").Value = 97
$ws.Range("K6").Value = 119.499996
$ws.Range("L6").Value = 291
$ws.Range("M6").Value = -6.499995999999996
$ws.Range("N6").Value = -517
$ws.Range("H23").Value = 2922.88
$ws.Range("J23").Value = 3635.3684
$ws.Range("L23").Value = 10906.1052
$ws.Range("N23").Value = -11376.1052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3124.75
$ws.Range("I43").Value = 3124.75
$ws.Range("K43").Value = 3124.75
$ws.Range("M43").Value = -2973.75
$ws.Range("H80").Value = 4147.778
$ws.Range("I80").Value = 1796.1111
$ws.Range("J80").Value = 6499.4443
$ws.Range("K80").Value = 1796.1111
$ws.Range("L80").Value = 6499.4443
$ws.Range("M80").Value = -798.1111000000001
$ws.Range("N80").Value = -8495.444299999999
$ws.Range("H83").Value = 4147.778
$ws.Range("I83").Value = 1796.1111
$ws.Range("J83").Value = 6499.4443
$ws.Range("K83").Value = 8980.5555
$ws.Range("L83").Value = 32497.2215
$ws.Range("M83").Value = -3988.5555
$ws.Range("N83").Value = -42481.2215
$ws.Range("H102").Value = 1337.8235
$ws.Range("I102").Value = 794.43634
$ws.Range("K102").Value = 794.43634
$ws.Range("M102").Value = 827.56366

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3923.25
$ws.Range("J7").Value = 4084.875
$ws.Range("L7").Value = 4084.875
$ws.Range("N7").Value = -4308.875
$ws.Range("H16").Value = 507.66666
$ws.Range("I16").Value = 217
$ws.Range("K16").Value = 217
$ws.Range("M16").Value = -47
$ws.Range("H46").Value = 2597.2307
$ws.Range("I46").Value = 2057
$ws.Range("J46").Value = 2796.2632
$ws.Range("K46").Value = 2057
$ws.Range("L46").Value = 2796.2632
$ws.Range("M46").Value = -1869
$ws.Range("N46").Value = -3172.2632
$ws.Range("H55").Value = 477.11765
$ws.Range("I55").Value = 227.33333
$ws.Range("J55").Value = 758.125
$ws.Range("K55").Value = 227.33333
$ws.Range("L55").Value = 758.125
$ws.Range("M55").Value = -54.33332999999999
$ws.Range("N55").Value = -1104.125
$ws.Range("H126").Value = 3923.25
$ws.Range("J126").Value = 4084.875
$ws.Range("L126").Value = 12254.625
$ws.Range("N126").Value = -17194.625
$ws.Range("H132").Value = 52634580
$ws.Range("I132").Value = 71431150
$ws.Range("J132").Value = 4172.6
$ws.Range("K132").Value = 214293450
$ws.Range("L132").Value = 12517.8
$ws.Range("M132").Value = -214290920
$ws.Range("N132").Value = -17577.8
$ws.Range("H136").Value = 4103.46
$ws.Range("I136").Value = 3721.738
$ws.Range("J136").Value = 6107.5
$ws.Range("K136").Value = 11165.214
$ws.Range("L136").Value = 18322.5
$ws.Range("M136").Value = -8615.214
$ws.Range("N136").Value = -23422.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 13893194
$ws.Range("I100").Value = 16671676
$ws.Range("K100").Value = 33343352
$ws.Range("M100").Value = -33342811
$ws.Range("H122").Value = 1951.742
$ws.Range("I122").Value = 1755.1364
$ws.Range("K122").Value = 5265.4092
$ws.Range("M122").Value = -2815.4092
$ws.Range("H132").Value = 5559.9707
$ws.Range("I132").Value = 6260.4165
$ws.Range("J132").Value = 3878.9
$ws.Range("K132").Value = 18781.2495
$ws.Range("L132").Value = 11636.7
$ws.Range("M132").Value = -16251.2495
$ws.Range("N132").Value = -16696.7
